$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.356.89"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "2.008.40"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'247.78"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("D7").Value = "'60.56"
$ws.Range("E7").Value = "  +4.61%  "
$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "'15.13"
$ws.Range("E12").Value = "  +10.38%  "
$ws.Range("D13").Value = "'22.80"
$ws.Range("E13").Value = "  +7.77%  "
$ws.Range("D14").Value = "'0.850"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").Value = "2.302.30"
$ws.Range("E15").Value = "  +2.94%  "
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "2.010.22"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").Value = "37.264.67"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "'70.59"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").Value = "'5.22"
$ws.Range("E21").Value = "  +4.01%  "
$ws.Range("D22").Value = "'231.09"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +7.32%  "
$ws.Range("D27").Value = "'9.46"
$ws.Range("E27").Value = "  +3.95%  "
$ws.Range("D28").Value = "'163.90"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "'19.73"
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("D30").Value = "'1.34"
$ws.Range("E30").Value = "  +15.55%  "
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "'4.84"
$ws.Range("E32").Value = "  +3.85%  "
$ws.Range("D33").Value = "'0.0655"
$ws.Range("E33").Value = "  +7.61%  "
$ws.Range("E34").Value = "  +4.99%  "
$ws.Range("E35").Value = "  +6.05%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").Value = "'3.29"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "'5.52"
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("D40").Value = "'0.0980"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").Value = "'16.71"
$ws.Range("E44").Value = "  +6.40%  "
$ws.Range("D45").Value = "'91.06"
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").Value = "1.374.55"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'7.28"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  +15.08%  "
$ws.Range("D51").Value = "'46.20"
$ws.Range("E51").Value = "  +5.79%  "
